$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated 2025 (column H) values per row, as described in the commit
# "Incorporo nuevos datos hasta diciembre de 2025"
$newValues = @{
    2  = 11.04662594326192
    3  = 11.73827866398485
    4  = 12.30585203826984
    5  = 11.78520396247183
    6  = 10.49487381709038
    7  = 11.35916506891029
    8  = 11.39276201862572
    9  = 10.92829416375072
    10 = 13.60046430866254
    11 = 11.50908501036198
    12 = 10.63149664638683
    13 = 11.51885833195096
    14 = 14.68713447755175
    15 = 11.2121250129479
    16 = 13.73060655955254
    17 = 14.91297330275808
    18 = 11.63063519186549
}

foreach ($row in $newValues.Keys) {
    $ws.Range("H$row").Value = $newValues[$row]
}
